$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = "'64.068.70"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = "'3.472.87"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = "'584.68"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = "'131.55"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -1.86%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = "'0.482"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").Value = "'7.65"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +5.60%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = "'0.123"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = "'0.386"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = "'4.066.48"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.120"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = "'0.0000177"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = "'3.475.36"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = "'64.067.52"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = "'24.84"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -3.66%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = "'10.00"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = "'5.68"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = "'13.40"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = "'384.63"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -2.26%  '
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = "'0.567"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("B23").Value = 'WrappedeETH'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D23").Value = "'3.614.62"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'74.67"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = "'5.59"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = "'0.0000111"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'2.22"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = "'7.10"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'7.94"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.43"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -4.52%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = "'0.153"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = "'3.502.36"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = "'22.93"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = "'5.22"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = "'6.76"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = "'1.50"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = "'161.93"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = "'0.0778"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = "'0.796"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = "'41.12"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = "'4.30"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = "'1.62"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'23.51"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -6.42%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = "'1.13"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'6.71"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = "'0.900"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = "'2.326.66"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -5.47%  '

Write-Output "applied cryptos update"
